$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing "2021" record): disambiguate id + fix dates
$ws.Cells.Item(2, 1).Value = "z0bug.li_partner_6_2021"
$ws.Cells.Item(2, 5).Value = "<###-01-25"
$ws.Cells.Item(2, 6).Value = "<###-01-25"
$ws.Cells.Item(2, 7).Value = "<###-01-25"
$ws.Cells.Item(2, 8).Value = "<###-12-31"

# Row 3 (existing "2022" record): disambiguate id + fix date
$ws.Cells.Item(3, 1).Value = "z0bug.li_partner_6_2022"
$ws.Cells.Item(3, 5).Value = "####-01-06"
$ws.Cells.Item(3, 6).Value = "####-01-06"
$ws.Cells.Item(3, 7).Value = "####-01-06"

# Selection moved to A4
$ws.Range("A4").Select()

# Column A widened to fit the longer id strings (closest reachable width
# given the engine's pixel-quantized column-width model; target OOXML
# width is 21.95, 21.1 "chars" is the nearest achievable setting)
$ws.Columns.Item(1).ColumnWidth = 21.1
